$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cluster name (col A) / Active cases (col B) for every data row, rows 1..101.
# Row 1 is the header (both cells are text); rows 2-101 pair a cluster name with a numeric count.
$rows = @(
    @(1, "Cluster name", "Active cases"),
    @(2, "3535 Opal Meadow Heights Aged Care Community Meadow Heights", 27),
    @(3, "Adorn Cosmetics Clayton", 5),
    @(4, "Al Haj Halal Meats Glenroy", 34),
    @(5, "Al-Taqwa College Truganina", 9),
    @(6, "Amiga Montessori Craigieburn", 28),
    @(7, "Anglicare Home Princes Hwy Dandenong", 5),
    @(8, "Best&Less Fountain Gate Narre Warren", 5),
    @(9, "Budget Car and Truck Rentals Campbellfield", 7),
    @(10, "CS Square Caroline Springs", 9),
    @(11, "Campbellfield Ford Complex Vaccination Clinic Campbellfield", 12),
    @(12, "Cannie Road Construction Site Cannie", 6),
    @(13, "Caroline Springs Police Station", 11),
    @(14, "Cedars Medical Clinic Coburg", 37),
    @(15, "Chemist Warehouse Campbellfield DC", 6),
    @(16, "Chemist Warehouse Fillo Drive Somerton", 11),
    @(17, "City of Moreland Community", 5),
    @(18, "City of Wyndham Community", 6),
    @(19, "Classy Cabinets and Kitchens Craigieburn", 9),
    @(20, "Coles Aurora Village Epping", 6),
    @(21, "Coles Broadmeadows Central Shopping Centre", 8),
    @(22, "Coles Campbellfield Plaza Campbellfield", 12),
    @(23, "Coles Coburg North Village", 29),
    @(24, "Coles Greenvale Shopping Centre", 5),
    @(25, "Coles Pakenham Place Shopping Centre", 13),
    @(26, "Coles Roxburgh Village Roxburgh Park", 8),
    @(27, "Community Kids Bayswater Early Education Centre Bayswater North", 8),
    @(28, "Community Kids Meadow Heights", 11),
    @(29, "Construction Site Olea Apartment Caulfield North", 16),
    @(30, "Costco Wholesale Epping", 27),
    @(31, "Crusader Caravans Epping", 14),
    @(32, "Crusader Caravans Epping", 22),
    @(33, "DRC Laverton Automotive Repairs Laverton North", 5),
    @(34, "Direct Freight Express Cambellfield", 13),
    @(35, "Don Watson Coldstore Derrimut", 5),
    @(36, "Epworth Healthcare Epworth Richmond Emergency Department", 6),
    @(37, "Fine Food Holdings Pty Ltd Dandenong South", 8),
    @(38, "Fitzroy Community School Fitzroy North", 35),
    @(39, "Fonterra Manufacturing Workplace Campbellfield", 9),
    @(40, "General Foods Campbellfield", 12),
    @(41, "Glenroy West Primary School", 6),
    @(42, "Goodstart Early Learning Altona", 9),
    @(43, "Green Leaves Early Learning Centre Highlands Craigieburn", 9),
    @(44, "Gumboots Early Learning Centre South Morang", 5),
    @(45, "Hamilton Marino 236 Jasper Road McKinnon", 13),
    @(46, "Health Care Providers Association South Melbourne", 7),
    @(47, "Hello Fresh Warehouse Ravenhall", 5),
    @(48, "IGA Meadow Heights Shopping Centre Meadow Heights", 6),
    @(49, "ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine", 10),
    @(50, "Ibis Kingsgate Hotel Melbourne", 6),
    @(51, "Ilim College Glenroy Campus Hadfield", 15),
    @(52, "Ilim Learning Sanctuary Glenroy", 10),
    @(53, "Industrial Galvanizers Valmont Coatings Campbellfield", 22),
    @(54, "KFC Fawkner", 6),
    @(55, "Kasr Sweets Coolaroo", 5),
    @(56, "Kids House Early Learning Cheltenham", 9),
    @(57, "Kippers Seafood Werribee", 6),
    @(58, "Level Crossing Removal Project Lilydale Construction Site John Street", 8),
    @(59, "Lineage Logistics Laverton North", 8),
    @(60, "Linfox Somerton National Distribution Centre Somerton", 9),
    @(61, "Mecca D.C Warehouse Melbourne Airport", 9),
    @(62, "Melbourne Assessment Prison West Melbourne", 5),
    @(63, "Melbourne Metropolitan Remand Centre Ravenhall", 11),
    @(64, "Melbourne Truck Repairs Campbellfield", 7),
    @(65, "Melbourne West Police Station Docklands", 7),
    @(66, "Mercy Hospital for Women Heidelberg", 5),
    @(67, "Mernda YMCA Early Learning Centre Mernda", 5),
    @(68, "Mill Park Police Station Mill Park", 5),
    @(69, "MyCentre Childcare Broadmeadows", 17),
    @(70, "National Gallery of Victoria Melbourne", 9),
    @(71, "Nido Early School Moonee Ponds", 14),
    @(72, "Northern Health Northern Hospital Epping Emergency Department Tier 1B", 52),
    @(73, "Northern Health The Northern Hospital Epping", 10),
    @(74, "OnQ Plumbing and Excavations Craigieburn", 18),
    @(75, "Oporto Coolaroo", 11),
    @(76, "Oscar Romero Catholic Primary School Craigieburn", 5),
    @(77, "Our Lady Help of Christian's Primary School Brunswick East", 10),
    @(78, "Paisley Park Early Learning Centre Bundoora", 6),
    @(79, "Panorama Construction Site Whitehorse Rd Box Hill", 14),
    @(80, "Private Residence Northern Community Services Fawkner", 5),
    @(81, "Ramsay Health Care Warringal Private Hospital Heidelberg", 9),
    @(82, "Richmond Quarter 261-271 Bridge Road Construction Site Richmond", 11),
    @(83, "Sacca's Fruit World Broadmeadows Central Shopping Centre", 6),
    @(84, "Salta Drive Construction Site Rangedale Drainage Altona North", 6),
    @(85, "Sharpline Stainless Steel Coburg North", 5),
    @(86, "St Margaret's Primary School OSHC Maribyrnong", 11),
    @(87, "St Vincents Hospital Emergency Department Melbourne", 6),
    @(88, "Tek Foods Somerton", 14),
    @(89, "The Homestead Child and Family Centre Roxburgh Park", 11),
    @(90, "The Huntly-Goornong Rail Works", 5),
    @(91, "The Royal Children's Hospital Melbourne Emergency Department Parkville Tier 1B", 16),
    @(92, "ThorwestenCabinets Pakenham", 13),
    @(93, "Total Window Concepts Hoppers Crossing", 6),
    @(94, "Unilodge College Square Student Accommodation 570 Lygon Street Carlton", 14),
    @(95, "Unilodge College Square Student Accommodation 570 Lygon Street Carlton", 14),
    @(96, "Werribee Mercy Hospital Emergency Department", 10),
    @(97, "Western Health Footscray Hospital Emergency Department", 6),
    @(98, "Western Health Sunshine Hospital Emergency Department", 10),
    @(99, "Woodlands Long Day Care and Kindergarten Roxburgh Park", 5),
    @(100, "Woolworths Greenvale Lakes Roxburgh Park", 5),
    @(101, "Yara Childcare Centre Truganina", 8)
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
}

